$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing row 6 (r2_adj row), shifting it down to row 7
$ws.Rows.Item(6).Insert()

# Set the new row's label, matching the bold/border/centered-top style used by
# the other label cells in column A (A2:A5, A7)
$ws.Range("A6").Value = "r2"
$ws.Range("A6").Borders.LineStyle = 1
$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").HorizontalAlignment = -4108
$ws.Range("A6").VerticalAlignment = -4160

# Fill in the new row's values
$ws.Range("B6").Value = 0.6
$ws.Range("C6").Value = 0.4
$ws.Range("D6").Value = 0.68
